# Remove the "All other diseases (Residual)" row (Excel row 4, rank 3) from the
# sheet. Every row below it (B = cause name, C = Z(10) value) moves up one
# position to fill the gap, while column A (the 1..N rank counter) is left
# untouched since it is independent of which row was dropped. The final row
# (136) becomes a duplicate of row 135 after the shift, so it is deleted
# outright, which also drops the sheet's row count/dimension back to 135 and
# lets the now-unreferenced shared string get cleaned up on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 4; $row -le 135; $row++) {
    $nextName  = $ws.Cells.Item($row + 1, 2).Value2
    $nextValue = $ws.Cells.Item($row + 1, 3).Value2
    $ws.Cells.Item($row, 2).Value = $nextName
    $ws.Cells.Item($row, 3).Value = $nextValue
}

$ws.Rows.Item(136).Delete()
